$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2417.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2417.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2417.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2555.5
$ws.Range("H106").Value = 2889.5
$ws.Range("I106").Value = 3990.3333
$ws.Range("J106").Value = 2417.7144
$ws.Range("K106").Value = 3990.3333
$ws.Range("L106").Value = 2417.7144
$ws.Range("M106").Value = -3359.3333
$ws.Range("N106").Value = -3679.7144
$ws.Range("H138").Value = 1877.9656
$ws.Range("I138").Value = 1507
$ws.Range("K138").Value = 4521
$ws.Range("M138").Value = 619

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3253.074
$ws.Range("I61").Value = 2672.111
$ws.Range("K61").Value = 2672.111
$ws.Range("M61").Value = -2460.111
$ws.Range("H74").Value = 1948.8966
$ws.Range("I74").Value = 2019.875
$ws.Range("J74").Value = 1861.5385
$ws.Range("K74").Value = 2019.875
$ws.Range("L74").Value = 1861.5385
$ws.Range("M74").Value = -1145.875
$ws.Range("N74").Value = -3609.5385
$ws.Range("H77").Value = 1948.8966
$ws.Range("I77").Value = 2019.875
$ws.Range("J77").Value = 1861.5385
$ws.Range("K77").Value = 10099.375
$ws.Range("L77").Value = 9307.692500000001
$ws.Range("M77").Value = -5731.375
$ws.Range("N77").Value = -18043.6925
$ws.Range("H136").Value = 3253.074
$ws.Range("I136").Value = 2672.111
$ws.Range("K136").Value = 8016.333
$ws.Range("M136").Value = -5466.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 2500
$ws.Range("I54").Value = 2500
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 2500
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -2016
$ws.Range("N54").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5651838
$ws.Range("I31").Value = 1494.0278
$ws.Range("J31").Value = 14495854
$ws.Range("K31").Value = 1494.0278
$ws.Range("L31").Value = 14495854
$ws.Range("M31").Value = -1199.0278
$ws.Range("N31").Value = -14496444
$ws.Range("H34").Value = 5651838
$ws.Range("I34").Value = 1494.0278
$ws.Range("J34").Value = 14495854
$ws.Range("K34").Value = 1494.0278
$ws.Range("L34").Value = 14495854
$ws.Range("M34").Value = -1292.0278
$ws.Range("N34").Value = -14496258
$ws.Range("H62").Value = 10064.4375
$ws.Range("I62").Value = 2585.9167
$ws.Range("J62").Value = 32500
$ws.Range("K62").Value = 2585.9167
$ws.Range("L62").Value = 32500
$ws.Range("M62").Value = -1961.9167
$ws.Range("N62").Value = -33748
$ws.Range("H65").Value = 10064.4375
$ws.Range("I65").Value = 2585.9167
$ws.Range("J65").Value = 32500
$ws.Range("K65").Value = 12929.5835
$ws.Range("L65").Value = 162500
$ws.Range("M65").Value = -9809.583500000001
$ws.Range("N65").Value = -168740
$ws.Range("H75").Value = 27400
$ws.Range("J75").Value = 27400
$ws.Range("L75").Value = 27400
$ws.Range("N75").Value = -29396
$ws.Range("H78").Value = 27400
$ws.Range("J78").Value = 27400
$ws.Range("L78").Value = 82200
$ws.Range("N78").Value = -92184
$ws.Range("H132").Value = 2835.76
$ws.Range("I132").Value = 1791.6
$ws.Range("J132").Value = 3531.8667
$ws.Range("K132").Value = 5374.799999999999
$ws.Range("L132").Value = 10595.6001
$ws.Range("M132").Value = -2844.799999999999
$ws.Range("N132").Value = -15655.6001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 67168.336
$ws.Range("I55").Value = 600
$ws.Range("J55").Value = 100452.5
$ws.Range("K55").Value = 1800
$ws.Range("L55").Value = 301357.5
$ws.Range("M55").Value = -1623
$ws.Range("N55").Value = -301711.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1208
$ws.Range("I22").Value = 245
$ws.Range("J22").Value = 1345.5714
$ws.Range("K22").Value = 245
$ws.Range("L22").Value = 1345.5714
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = -1935.5714
$ws.Range("H27").Value = 1208
$ws.Range("I27").Value = 245
$ws.Range("J27").Value = 1345.5714
$ws.Range("K27").Value = 245
$ws.Range("L27").Value = 1345.5714
$ws.Range("M27").Value = -138
$ws.Range("N27").Value = -1559.5714
$ws.Range("H41").Value = 1000
$ws.Range("I41").Value = 1000
$ws.Range("K41").Value = 1000
$ws.Range("M41").Value = -562
$ws.Range("H46").Value = 477085.28
$ws.Range("I46").Value = 777.0769
$ws.Range("J46").Value = 1251086.1
$ws.Range("K46").Value = 777.0769
$ws.Range("L46").Value = 1251086.1
$ws.Range("M46").Value = -589.0769
$ws.Range("N46").Value = -1251462.1
$ws.Range("H55").Value = 245.38889
$ws.Range("I55").Value = 152
$ws.Range("J55").Value = 376.13333
$ws.Range("K55").Value = 152
$ws.Range("L55").Value = 376.13333
$ws.Range("M55").Value = 21
$ws.Range("N55").Value = -722.13333
$ws.Range("H136").Value = 29413058
$ws.Range("I136").Value = 33334460
$ws.Range("J136").Value = 2542.5
$ws.Range("K136").Value = 100003380
$ws.Range("L136").Value = 7627.5
$ws.Range("M136").Value = -100000830
$ws.Range("N136").Value = -12727.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 3000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3000
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -3586
$ws.Range("H62").Value = 3696.55
$ws.Range("I62").Value = 2920.8333
$ws.Range("J62").Value = 4860.125
$ws.Range("K62").Value = 2920.8333
$ws.Range("L62").Value = 4860.125
$ws.Range("M62").Value = -2296.8333
$ws.Range("N62").Value = -6108.125
$ws.Range("H65").Value = 3696.55
$ws.Range("I65").Value = 2920.8333
$ws.Range("J65").Value = 4860.125
$ws.Range("K65").Value = 14604.1665
$ws.Range("L65").Value = 24300.625
$ws.Range("M65").Value = -11484.1665
$ws.Range("N65").Value = -30540.625
$ws.Range("H70").Value = 20105
$ws.Range("J70").Value = 20105
$ws.Range("L70").Value = 20105
$ws.Range("N70").Value = -20735
$ws.Range("H73").Value = 20105
$ws.Range("J73").Value = 20105
$ws.Range("L73").Value = 20105
$ws.Range("N73").Value = -22289
